$d = $word.ActiveDocument

function SplitRunAt($absOffset, $tempName) {
    # Forces a run boundary at a given absolute character offset by adding
    # then immediately deleting a zero-length bookmark there. Word (and this
    # COM host) always places bookmarkStart/End on a run boundary, so the
    # surrounding text gets split into separate <w:r> runs that persist
    # after the temporary bookmark is removed.
    $bm = $d.Range($absOffset, $absOffset)
    $d.Bookmarks.Add($tempName, $bm) | Out-Null
    $d.Bookmarks.Item($tempName).Delete()
}

# =====================================================================
# Part 1: "OUTPUTS ... shoot" paragraph -> "OUTPUTS ... velocity", with
# the _GoBack bookmark relocated into the middle of this paragraph.
# =====================================================================

$outputsIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "OUTPUTS*shoot*") {
        $outputsIdx = $i
        break
    }
}

$pOut = $d.Paragraphs.Item($outputsIdx)
$outStart = $pOut.Range.Start
$outEnd = $pOut.Range.End

# "OUTPUTS     turn (0-left 0.5-straight 1-right) " == 47 chars (kept as-is)
$prefixLen = 47

# Replace the trailing "     shoot  " with " velocity  "
$tailRange = $d.Range($outStart + $prefixLen, $outEnd)
$tailRange.Text = " velocity  "

$pOut = $d.Paragraphs.Item($outputsIdx)
$outStart = $pOut.Range.Start

# Split so "velocity" and the trailing "  " become separate runs.
SplitRunAt ($outStart + $prefixLen + 9) "TmpSplitOut1"
# Split so the lone extra space becomes its own run (before the bookmark).
SplitRunAt ($outStart + $prefixLen) "TmpSplitOut2"

# Move _GoBack here: delete it from its old spot, re-add between the lone
# space run and "velocity".
$d.Bookmarks.Item("_GoBack").Delete()
$bmPos = $outStart + $prefixLen + 1
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos)) | Out-Null

# =====================================================================
# Part 2: "- calculate_fitness" / "- " paragraphs -> "- calculate_fitness
# -> " and "- make_decision -> output units".
# =====================================================================

$calcIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq "- calculate_fitness") {
        $calcIdx = $i
        break
    }
}
$blankIdx = $calcIdx + 1

# --- "- calculate_fitness" -> "- calculate_fitness -> " (4 runs) ---
$pCalc = $d.Paragraphs.Item($calcIdx)
$calcStart = $pCalc.Range.Start
$calcTextLen = "- calculate_fitness".Length

$appendRange = $d.Range($calcStart + $calcTextLen, $calcStart + $calcTextLen)
$appendRange.Text = " -> "

$pCalc = $d.Paragraphs.Item($calcIdx)
$calcStart = $pCalc.Range.Start

SplitRunAt ($calcStart + 19) "TmpSplitCalc3"   # before " -> "
SplitRunAt ($calcStart + 11) "TmpSplitCalc2"   # before "_fitness"
SplitRunAt ($calcStart + 2)  "TmpSplitCalc1"   # before "calculate"

# --- "- " (bookmarked) -> "- make_decision -> output units" (2 runs) ---
$pBlank = $d.Paragraphs.Item($blankIdx)
$blankStart = $pBlank.Range.Start

$fillRange = $d.Range($blankStart + 2, $blankStart + 2)
$fillRange.Text = "make_decision -> output units"

$pBlank = $d.Paragraphs.Item($blankIdx)
$blankStart = $pBlank.Range.Start
SplitRunAt ($blankStart + 2) "TmpSplitMake1"   # before "make_decision -> output units"

# The old _GoBack bookmark that used to live in this paragraph was already
# removed in Part 1 (deleting by name relocates it, it cannot exist twice).

Write-Output "done"
